$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
$oldHeaders = @(
  "Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old",
  "Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old"
)
$newHeadersFV2404 = @(
  "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
  "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404"
)
for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $newHeadersFV2404[$i]
}

$newHeaders = @(
  "Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new",
  "Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new"
)
$newHeadersFV2410 = @(
  "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
  "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410"
)
# columns L..U are indices 12..21 (diff is column K / index 11, unchanged)
for ($i = 0; $i -lt $newHeadersFV2410.Length; $i++) {
  $ws.Cells.Item(1, 11 + $i + 1).Value = $newHeadersFV2410[$i]
}

# 2. Freeze the header row (top row frozen, matching pane ySplit=1 / topLeftCell A2)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn the used range into a table (ListObject) with autofilter, matching xl/tables/table1.xml
$rng = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
